$wb = $excel.ActiveWorkbook

# 1. Rename sheet "DataListParsing" -> "DataCollectionParsing"
$ws4 = $wb.Worksheets.Item("DataListParsing")
$ws4.Name = "DataCollectionParsing"

# 2. Update the "dataList" value cells to "dataCollection"
$ws4.Range("B1").Value = "dataCollection"

$ws5 = $wb.Worksheets.Item("RowValueListParsing")
$ws5.Range("B1").Value = "dataCollection"

$wsLists = $wb.Worksheets.Item(".lists")
$wsLists.Range("E3").Value = "dataCollection"
